$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 744.8333
$ws.Range("I28").Value = 514.7857
$ws.Range("J28").Value = 1550
$ws.Range("K28").Value = 514.7857
$ws.Range("L28").Value = 1550
$ws.Range("M28").Value = -29.78570000000002
$ws.Range("N28").Value = -2520
$ws.Range("H32").Value = 983.1667
$ws.Range("J32").Value = 974.75
$ws.Range("L32").Value = 974.75
$ws.Range("N32").Value = -1626.75
$ws.Range("H33").Value = 407.7931
$ws.Range("I33").Value = 369.16666
$ws.Range("K33").Value = 369.16666
$ws.Range("M33").Value = -140.16666
$ws.Range("H40").Value = 2622.2856
$ws.Range("I40").Value = 4925
$ws.Range("K40").Value = 4925
$ws.Range("M40").Value = -4750
$ws.Range("H137").Value = 3233.9783
$ws.Range("I137").Value = 1666.5667
$ws.Range("J137").Value = 6172.875
$ws.Range("K137").Value = 4999.7001
$ws.Range("L137").Value = 18518.625
$ws.Range("M137").Value = -2449.7001
$ws.Range("N137").Value = -23618.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6367.6343
$ws.Range("I32").Value = 5742.169
$ws.Range("K32").Value = 5742.169
$ws.Range("M32").Value = -5455.169
$ws.Range("H88").Value = 4681.7646
$ws.Range("I88").Value = 14178
$ws.Range("J88").Value = 1759.8462
$ws.Range("K88").Value = 14178
$ws.Range("L88").Value = 1759.8462
$ws.Range("M88").Value = -13772
$ws.Range("N88").Value = -2571.8462
$ws.Range("H91").Value = 4681.7646
$ws.Range("I91").Value = 14178
$ws.Range("J91").Value = 1759.8462
$ws.Range("K91").Value = 14178
$ws.Range("L91").Value = 1759.8462
$ws.Range("M91").Value = -12774
$ws.Range("N91").Value = -4567.8462
$ws.Range("H132").Value = 4969.5957
$ws.Range("I132").Value = 1733
$ws.Range("J132").Value = 8647.546
$ws.Range("K132").Value = 5199
$ws.Range("L132").Value = 25942.638
$ws.Range("M132").Value = -2669
$ws.Range("N132").Value = -31002.638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4884.515
$ws.Range("I105").Value = 4507.2
$ws.Range("J105").Value = 6063.625
$ws.Range("K105").Value = 4507.2
$ws.Range("L105").Value = 6063.625
$ws.Range("M105").Value = -2760.2
$ws.Range("N105").Value = -9557.625
$ws.Range("H107").Value = 3334.5557
$ws.Range("I107").Value = 3002.75
$ws.Range("J107").Value = 3600
$ws.Range("K107").Value = 3002.75
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = -1082.75
$ws.Range("N107").Value = -7440
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1494.32
$ws.Range("I31").Value = 1090.3944
$ws.Range("J31").Value = 2483.2415
$ws.Range("K31").Value = 1090.3944
$ws.Range("L31").Value = 2483.2415
$ws.Range("M31").Value = -795.3943999999999
$ws.Range("N31").Value = -3073.2415
$ws.Range("H34").Value = 1494.32
$ws.Range("I34").Value = 1090.3944
$ws.Range("J34").Value = 2483.2415
$ws.Range("K34").Value = 1090.3944
$ws.Range("L34").Value = 2483.2415
$ws.Range("M34").Value = -888.3943999999999
$ws.Range("N34").Value = -2887.2415
$ws.Range("H105").Value = 334.33334
$ws.Range("I105").Value = 334.33334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 334.33334
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1412.66666
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2177.0784
$ws.Range("I132").Value = 2009.1818
$ws.Range("J132").Value = 2304.4482
$ws.Range("K132").Value = 6027.5454
$ws.Range("L132").Value = 6913.344599999999
$ws.Range("M132").Value = -3497.5454
$ws.Range("N132").Value = -11973.3446
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 217.9
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 217.9
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 653.7
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -991.7
$ws.Range("H56").Value = 103379.9
$ws.Range("I56").Value = 103379.9
$ws.Range("K56").Value = 103379.9
$ws.Range("M56").Value = -102849.9
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H103").Value = 3514
$ws.Range("I103").Value = 5000
$ws.Range("J103").Value = 2028
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 6084
$ws.Range("M103").Value = -14121
$ws.Range("N103").Value = -7842
$ws.Range("H118").Value = 8430.799999999999
$ws.Range("J118").Value = 10031.25
$ws.Range("L118").Value = 30093.75
$ws.Range("N118").Value = -32579.75
$ws.Range("H121").Value = 1461.3529
$ws.Range("I121").Value = 601.6667
$ws.Range("J121").Value = 1930.2727
$ws.Range("K121").Value = 1805.0001
$ws.Range("L121").Value = 5790.8181
$ws.Range("M121").Value = -495.0001
$ws.Range("N121").Value = -8410.8181
$ws.Range("H122").Value = 697.0645
$ws.Range("I122").Value = 385.65
$ws.Range("J122").Value = 1263.2727
$ws.Range("K122").Value = 3470.85
$ws.Range("L122").Value = 11369.4543
$ws.Range("M122").Value = -1020.85
$ws.Range("N122").Value = -16269.4543
$ws.Range("H131").Value = 12587
$ws.Range("I131").Value = 536.5179000000001
$ws.Range("J131").Value = 40704.793
$ws.Range("K131").Value = 1609.5537
$ws.Range("L131").Value = 122114.379
$ws.Range("M131").Value = 3430.4463
$ws.Range("N131").Value = -132194.379
$ws.Range("H132").Value = 1791.4073
$ws.Range("I132").Value = 1703.1666
$ws.Range("J132").Value = 1862
$ws.Range("K132").Value = 15328.4994
$ws.Range("L132").Value = 16758
$ws.Range("M132").Value = -12798.4994
$ws.Range("N132").Value = -21818
$ws.Range("H137").Value = 26346056
$ws.Range("I137").Value = 41669044
$ws.Range("J137").Value = 78076.14
$ws.Range("K137").Value = 125007132
$ws.Range("L137").Value = 234228.42
$ws.Range("M137").Value = -125002032
$ws.Range("N137").Value = -244428.42
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1107.1428
$ws.Range("I43").Value = 1107.1428
$ws.Range("K43").Value = 1107.1428
$ws.Range("M43").Value = -956.1428000000001
$ws.Range("H80").Value = 6581.231
$ws.Range("I80").Value = 14187.5
$ws.Range("J80").Value = 3200.6667
$ws.Range("K80").Value = 14187.5
$ws.Range("L80").Value = 3200.6667
$ws.Range("M80").Value = -13189.5
$ws.Range("N80").Value = -5196.6667
$ws.Range("H83").Value = 6581.231
$ws.Range("I83").Value = 14187.5
$ws.Range("J83").Value = 3200.6667
$ws.Range("K83").Value = 70937.5
$ws.Range("L83").Value = 16003.3335
$ws.Range("M83").Value = -65945.5
$ws.Range("N83").Value = -25987.3335
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 2475.6667
$ws.Range("I113").Value = 2811.7144
$ws.Range("J113").Value = 1299.5
$ws.Range("K113").Value = 2811.7144
$ws.Range("L113").Value = 1299.5
$ws.Range("M113").Value = -641.7143999999998
$ws.Range("N113").Value = -5639.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 39000
$ws.Range("J6").Value = 39000
$ws.Range("L6").Value = 39000
$ws.Range("N6").Value = -39224
$ws.Range("H68").Value = 4371.4287
$ws.Range("I68").Value = 4640
$ws.Range("J68").Value = 3700
$ws.Range("K68").Value = 4640
$ws.Range("L68").Value = 3700
$ws.Range("M68").Value = -3891
$ws.Range("N68").Value = -5198
$ws.Range("H71").Value = 4371.4287
$ws.Range("I71").Value = 4640
$ws.Range("J71").Value = 3700
$ws.Range("K71").Value = 23200
$ws.Range("L71").Value = 18500
$ws.Range("M71").Value = -19456
$ws.Range("N71").Value = -25988
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 57433.332
$ws.Range("I109").Value = 40000
$ws.Range("J109").Value = 60920
$ws.Range("K109").Value = 40000
$ws.Range("L109").Value = 60920
$ws.Range("M109").Value = -38613
$ws.Range("N109").Value = -63694
